# Add a new LeetCode "Top SQL 50" entry (1789. Primary Department for Each
# Employee) as a new row appended to the Table2 table on Sheet1, matching
# the formatting conventions already used by the sheet (green fill for
# "Easy" difficulty, Hyperlink style + live hyperlink for the Link column),
# and move the active selection the way the author's session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the worksheet table by one row; this extends Table2's ref/autoFilter
# from A1:E26 to A1:E27 (mirrors the table1.xml change in the diff).
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

# New row index.
$r = 27

# Match the green "Easy" fill used on every other Easy row (B2, B3, ...)
# before writing the value, so the engine reuses the existing style index
# instead of minting a fresh one.
$ws.Range("B$r").Interior.Color = 5287936

$ws.Range("A$r").Value = "1789. Primary Department for Each Employee"
$ws.Range("B$r").Value = "Easy"
$ws.Range("C$r").Value = "Advanced Select and Joins"
$ws.Range("D$r").Value = "Use union function"

# Link column: set the display text, wire up a real hyperlink relationship,
# then re-apply the built-in "Hyperlink" style (Hyperlinks.Add already
# applies hyperlink formatting, but re-asserting the named style keeps the
# cell on the same style index as the rest of column E).
$linkUrl = "https://leetcode.com/problems/primary-department-for-each-employee/solutions/3871250/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 "
$ws.Range("E$r").Value = $linkUrl
$null = $ws.Hyperlinks.Add($ws.Range("E$r"), $linkUrl)
$ws.Range("E$r").Style = "Hyperlink"

# Restore the selection state recorded in the saved workbook.
$ws.Range("D31").Select()
